$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet 1")
$ws2 = $wb.Worksheets.Item("References")

$ws1.Range("E5").Value = "Based on official disease reports to the WOAH"
$ws1.Range("E6").Value = "Venezuelan equine encephalitis is a disease listed in the World Organisation for Animal Health ({ref009:WOAH}) Terrestrial Animal Health Code. The map to the right displays occurrence reported to the {ref001:WOAH-WAHIS} system since 2005."
$ws1.Range("E7").Value = "As described in the WOAH {ref005:Terrestrial Animal Health Code}, the WOAH early warning system includes immediate notifications and follow-up reports on:"
$ws1.Range("E14").Value = "Information on stable situations (disease present or absent in a zone or country) is provided by countries through the WOAH monitoring system, which is a different reporting channel. This information is available in a different spatial and temporal scale, which can be browsed on the map independently from the outbreak notification points."
$ws1.Range("E17").Value = "For more up to date reports, visit the original data source: {ref001:WOAH-WAHIS}."
$ws1.Range("E72").Value = "WOAH-prescribed tests for international trade include ({ref010:WOAH, Terrestrial Manual}):"
$ws1.Range("E140").Value = "Geographical distribution data has been kindly provided by the World Organisation of Animal Health (WOAH). {ref001:WOAH-WAHIS} (WOAH World Animal Health Information System) is the original source of these data."

$ws2.Range("C2").Value = "WOAH-WAHIS (WOAH World Animal Health Information System)"
$ws2.Range("C5").Value = "WOAH (World Organisation for Animal Health). Terrestrial Animal Health Code 2021. WOAH, Paris, France"
$ws2.Range("C8").Value = "WOAH (World Organisation for Animal Health), 2013. VENEZUELAN EQUINE ENCEPHALITIS. Technical disease card. "
$ws2.Range("C9").Value = "WOAH (World Organisation for Animal Health), 2021. VENEZUELAN EQUINE ENCEPHALOMYELITIS . Chapter 12.11. WOAH Terrestrial Animal Health Code, Paris, France"
$ws2.Range("C10").Value = "WOAH (World Organisation for Animal Health), 2019. EQUINE ENCEPHALOMYELITIS (Eastern, Western and Venezuelan). Chapter 3.5.5. WOAH Terrestrial Manual, Paris, France"
